# Updates cryptos price (D) and volume-change (E) columns to match the
# latest scrape, per the GitHub Actions commit on 2023-08-04.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.158.22"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.831.22"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'241.61"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'0.6596"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.07402"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.2927"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "'22.91"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'0.07754"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.821.34"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "'4.993"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'0.6656"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'82.75"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D17").Value = "'0.000008423"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("D18").Value = "29.143.50"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "2.075.59"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'226.97"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.124"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'158.76"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'8.603"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "'17.90"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "'1.519"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").Value = "'4.114"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").Value = "'4.042"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").Value = "'1.188"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "'0.05253"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").Value = "'1.861"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").Value = "'0.7405"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "'1.142"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "'2.653"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "1.302.22"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "'0.01790"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "'2.731"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").Value = "'0.9206"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'5.946"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "'0.08478"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("D45").Value = "'102.12"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "1.976.15"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.5138"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  -9.53%  "
$ws.Range("D49").Value = "'1.750"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'63.28"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'0.05842"
$ws.Range("E51").Value = "  -1.40%  "
